$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Add the two new sheets: "QuickLinks" and "AddRelationship".
# They must end up positioned right after "Users" and before
# "RecentlyViewedListView".
# ---------------------------------------------------------------------------
$usersSheet = $wb.Worksheets.Item("Users")
$recentSheet = $wb.Worksheets.Item("RecentlyViewedListView")

# "QuickLinks" is created first so it receives the lower internal
# sheetId (8); "AddRelationship" is created second (sheetId 9). Both are
# inserted immediately before RecentlyViewedListView.
$quickLinks = $wb.Worksheets.Add($recentSheet)
$quickLinks.Name = "QuickLinks"

$addRelationship = $wb.Worksheets.Add($recentSheet)
$addRelationship.Name = "AddRelationship"

# Re-fetch by name (handles returned from Add() can go stale once the
# worksheets collection changes again) and reorder so QuickLinks sits
# ahead of AddRelationship: Users, QuickLinks, AddRelationship,
# RecentlyViewedListView, Activity, Contact, ContactTypes.
$quickLinks = $wb.Worksheets.Item("QuickLinks")
$addRelationship = $wb.Worksheets.Item("AddRelationship")
$quickLinks.Move($addRelationship)

$quickLinks = $wb.Worksheets.Item("QuickLinks")
$addRelationship = $wb.Worksheets.Item("AddRelationship")

# ---------------------------------------------------------------------------
# Populate "QuickLinks" - single column list of quick link names.
# Row 5 ("Engagements Contacts") is written last (after rows 6-15) so the
# shared-string table gets new entries in the same order as the source
# workbook.
# ---------------------------------------------------------------------------
$quickLinks.Cells.Item(1, 1).Value = "QuickLinkNames"
$quickLinks.Cells.Item(2, 1).Value = "HL Relationships"
$quickLinks.Cells.Item(3, 1).Value = "Industry Focus"
$quickLinks.Cells.Item(4, 1).Value = "Opportunity Contacts"
$quickLinks.Cells.Item(6, 1).Value = "Engagements Shown"
$quickLinks.Cells.Item(7, 1).Value = "Affiliated Companies"
$quickLinks.Cells.Item(8, 1).Value = "Memberships"
$quickLinks.Cells.Item(9, 1).Value = "Contact Sectors"
$quickLinks.Cells.Item(10, 1).Value = "Campaign History"
$quickLinks.Cells.Item(11, 1).Value = "Contact Email History"
$quickLinks.Cells.Item(12, 1).Value = "Contact Sources"
$quickLinks.Cells.Item(13, 1).Value = "Development Leads"
$quickLinks.Cells.Item(14, 1).Value = "Files"
$quickLinks.Cells.Item(15, 1).Value = "Contact History"
$quickLinks.Cells.Item(5, 1).Value = "Engagements Contacts"

$quickLinks.Range("A1").Font.Bold = $true
$quickLinks.Columns.Item(1).AutoFit() | Out-Null
$quickLinks.Activate()
$quickLinks.Range("F14").Select() | Out-Null

# ---------------------------------------------------------------------------
# Populate "AddRelationship" - header row + one data row. Column B of the
# header is written before column A so the shared strings come out in the
# same order as the source workbook ("Strength Rating" before "Lookup
# Employee").
# ---------------------------------------------------------------------------
$addRelationship.Cells.Item(1, 2).Value = "Strength Rating"
$addRelationship.Cells.Item(1, 1).Value = "Lookup Employee"
$addRelationship.Cells.Item(1, 3).Value = "Type"
$addRelationship.Cells.Item(1, 4).Value = "Personal Note"
$addRelationship.Cells.Item(1, 5).Value = "Outlook Categories"

$addRelationship.Cells.Item(2, 1).Value = "James Craven"
$addRelationship.Cells.Item(2, 2).Value = "High"
$addRelationship.Cells.Item(2, 3).Value = "Business"
$addRelationship.Cells.Item(2, 4).Value = "Test"
$addRelationship.Cells.Item(2, 5).Value = "Test"

$addRelationship.Range("A1:E1").Font.Bold = $true
$addRelationship.Columns.Item(1).AutoFit() | Out-Null
$addRelationship.Columns.Item(2).AutoFit() | Out-Null
$addRelationship.Columns.Item(4).AutoFit() | Out-Null
$addRelationship.Columns.Item(5).AutoFit() | Out-Null

# ---------------------------------------------------------------------------
# "AddRelationship" is the newly-active/selected tab.
# ---------------------------------------------------------------------------
$addRelationship.Range("E2").Select() | Out-Null
$addRelationship.Activate()
